# Manuscript revision: split the combined "Etoh_and_smoking" model into
# separate Male / Female models.
#
# - Sheet 2 "Model Estimates" is renamed to "Parameter Estimates".
# - Sheet 1 "Model Summaries" row 2 is updated with the Male-model summary
#   stats, and a new row 3 is added with the Female-model summary stats.
# - Sheet 2 term-level estimates are rewritten: rows 2-15 become the Male
#   model (the old "sexMale" term is dropped since the models are now
#   sex-stratified), and rows 16-28 are added for the Female model.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws2.Name = "Parameter Estimates"

# ---------------------------------------------------------------------
# Sheet 1: "Model Summaries"
# ---------------------------------------------------------------------

# Row 2 -> Males (replaces the old combined "Etoh_and_smoking" row)
$ws1.Cells.Item(2, 1).Value = "Males"
$ws1.Cells.Item(2, 2).Value = 0.0533119027056901
$ws1.Cells.Item(2, 3).Value = 0.0444579852489807
$ws1.Cells.Item(2, 4).Value = 15.4326147813358
$ws1.Cells.Item(2, 5).Value = 6.02127848676646
$ws1.Cells.Item(2, 6).Value = 0.0000000000510341712493853
$ws1.Cells.Item(2, 7).Value = 13
$ws1.Cells.Item(2, 8).Value = -5827.17685708492
$ws1.Cells.Item(2, 9).Value = 11684.3537141698
$ws1.Cells.Item(2, 10).Value = 11763.0599229386
$ws1.Cells.Item(2, 11).Value = 331050.182594856
$ws1.Cells.Item(2, 12).Value = 1390
$ws1.Cells.Item(2, 13).Value = 1404

# Row 3 -> Females (new row)
$ws1.Cells.Item(3, 1).Value = "Females"
$ws1.Cells.Item(3, 2).Value = 0.0638559977684699
$ws1.Cells.Item(3, 3).Value = 0.0568741904180795
$ws1.Cells.Item(3, 4).Value = 14.3070124078752
$ws1.Cells.Item(3, 5).Value = 9.14605553602079
$ws1.Cells.Item(3, 6).Value = 0.0000000000000000293971164399593
$ws1.Cells.Item(3, 7).Value = 12
$ws1.Cells.Item(3, 8).Value = -6610.72828063101
$ws1.Cells.Item(3, 9).Value = 13249.456561262
$ws1.Cells.Item(3, 10).Value = 13324.9363745475
$ws1.Cells.Item(3, 11).Value = 329347.181898902
$ws1.Cells.Item(3, 12).Value = 1609
$ws1.Cells.Item(3, 13).Value = 1622

# ---------------------------------------------------------------------
# Sheet 2: "Parameter Estimates"
# ---------------------------------------------------------------------

# Male-model terms (rows 2-15); the "sexMale" term from the old combined
# model is gone now that the models are sex-stratified.
$maleRows = @(
    @("(Intercept)", 5.29333136071519, 16.809180251152, 0.314907168679591, 0.752879403039063, -27.6807688356121, 38.2674315570424),
    @("etoh_g_day", 0.113237786158915, 0.0178410797894489, 6.34702537600236, 0.000000000296465376974968, 0.0782394374115014, 0.148236134906329),
    @("smok_cat5Ex-smoker who smoked occasionally", 3.47477401755889, 1.59045777468678, 2.18476345166926, 0.0290727337511521, 0.354817348823949, 6.59473068629383),
    @("smok_cat5Ex-smoker who smoked regularly", 0.747932910436062, 0.961957574544779, 0.77751132713935, 0.436989461344672, -1.13911244059951, 2.63497826147164),
    @("smok_cat5Occasional smoker", 4.4398065829523, 2.91271718463475, 1.52428344446666, 0.127665395262149, -1.2739895007757, 10.1536026666803),
    @("smok_cat5Regular smoker", 4.68426062357128, 1.40476429439318, 3.33455273761408, 0.000876872747098563, 1.92857367730459, 7.43994756983796),
    @("age_at_survey", 0.276317934426404, 0.0767977202162969, 3.59799657656723, 0.000331928287259142, 0.125665988083683, 0.426969880769125),
    @("bmi", -0.0953909786765939, 0.0929649142098184, -1.02609656005598, 0.305024600151243, -0.277757658675009, 0.0869757013218209),
    @("educ_levelSecondary school", 2.42666003318875, 15.9596554273863, 0.15204965071017, 0.879169840929234, -28.8809510042636, 33.7342710706411),
    @("educ_levelVocational school / apprenticeship", 5.66355371981544, 15.5336572744265, 0.36459885909415, 0.715466326014674, -24.8083886193858, 36.1354960590167),
    @("educ_levelSecondary school / Middle school", 4.29614047481476, 15.5660212081408, 0.275994772033842, 0.782593108311736, -26.2392892907272, 34.8315702403567),
    @("educ_levelTechnical college / technician / master", 7.5602729996053, 15.5521223824365, 0.486124839664545, 0.626955282996558, -22.9478918271034, 38.068437826314),
    @("educ_levelAbitur", 4.20148907317058, 15.6094242954748, 0.269163615110942, 0.787843725329472, -26.4190833185554, 34.8220614648966),
    @("educ_levelUniversity", 5.11101407623988, 15.5589641512687, 0.328493209865974, 0.742588256455526, -25.4105720576131, 35.6326002100929)
)

$r = 2
foreach ($row in $maleRows) {
    $ws2.Cells.Item($r, 1).Value = "Males"
    $ws2.Cells.Item($r, 2).Value = $row[0]
    $ws2.Cells.Item($r, 3).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]
    $ws2.Cells.Item($r, 5).Value = $row[3]
    $ws2.Cells.Item($r, 6).Value = $row[4]
    $ws2.Cells.Item($r, 7).Value = $row[5]
    $ws2.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# Female-model terms (rows 16-28); note "educ_levelSecondary school" is not
# a term in the Female model.
$femaleRows = @(
    @("(Intercept)", 30.4499099783251, 5.10902068914304, 5.96002870824792, 0.00000000309267038358885, 20.4288752398763, 40.4709447167738),
    @("etoh_g_day", 0.228322918891214, 0.0293875122668693, 7.76938574513567, 0.0000000000000139550232729419, 0.170681092948191, 0.285964744834237),
    @("smok_cat5Ex-smoker who smoked occasionally", -0.248597156307385, 1.25793878202011, -0.197622618732022, 0.843365298871899, -2.71596791099108, 2.21877359837631),
    @("smok_cat5Ex-smoker who smoked regularly", 1.18537578073293, 0.862620198596931, 1.37415722778225, 0.169584255971074, -0.506601508495347, 2.8773530699612),
    @("smok_cat5Occasional smoker", 2.25780992833913, 3.0321148037793, 0.744632071821602, 0.45660291499987, -3.68949966940958, 8.20511952608784),
    @("smok_cat5Regular smoker", 1.33412879987719, 1.12849957751242, 1.1822147092142, 0.237295195726541, -0.879354792508407, 3.54761239226279),
    @("age_at_survey", -0.0192072251295009, 0.0669066469133445, -0.287074992031473, 0.774091864217458, -0.150440561955978, 0.112026111696977),
    @("bmi", -0.251611075332022, 0.0627095545498564, -4.01232439200286, 0.0000628834611220802, -0.374612069614973, -0.128610081049072),
    @("educ_levelVocational school / apprenticeship", 2.87583678167169, 1.78043571608051, 1.61524325517498, 0.10645408182697, -0.616380072789792, 6.36805363613316),
    @("educ_levelSecondary school / Middle school", 4.06403313488797, 1.84711680898249, 2.20020364447157, 0.0279341484730318, 0.441025354136364, 7.68704091563958),
    @("educ_levelTechnical college / technician / master", 3.90103786013192, 1.95866602164811, 1.99168097930724, 0.0465748687994227, 0.0592330528691005, 7.74284266739474),
    @("educ_levelAbitur", 2.83329071021574, 2.24932236259308, 1.25961967805693, 0.207989274096524, -1.57861891101848, 7.24520033144997),
    @("educ_levelUniversity", 4.19350854700612, 1.97798064642134, 2.12009584350242, 0.0341506541422248, 0.313819272721278, 8.07319782129096)
)

foreach ($row in $femaleRows) {
    $ws2.Cells.Item($r, 1).Value = "Females"
    $ws2.Cells.Item($r, 2).Value = $row[0]
    $ws2.Cells.Item($r, 3).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]
    $ws2.Cells.Item($r, 5).Value = $row[3]
    $ws2.Cells.Item($r, 6).Value = $row[4]
    $ws2.Cells.Item($r, 7).Value = $row[5]
    $ws2.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}
